# Commit: "criar tutorial e deck"
# Recolor two bullet items (paragraph mark + run text) to RGB #81D41A.
#   - "Pagina de criação de tutorial, deck e notícias;" goes from the
#     default (no explicit color) to #81D41A.
#   - "Adicionar hover e active nos botões e links" goes from explicit
#     "auto" to #81D41A.
#
# wdColor value for #81D41A (R=0x81=129, G=0xD4=212, B=0x1A=26) using
# Word's RGB(r,g,b) = r + g*256 + b*65536 packing:
$greenColor = 129 + (212 * 256) + (26 * 65536)   # 1758337

$d = $word.ActiveDocument
$paras = $d.Paragraphs

for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text

    if ($t -like "*Pagina de cria*tutorial, deck e not*cias*") {
        $p.Range.Font.Color = $greenColor
    }
    elseif ($t -like "*Adicionar hover e active nos bot*") {
        $p.Range.Font.Color = $greenColor
    }
}
